$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data updates per diff
# J2 holds a zero-padded code ("001"); force text format so Excel keeps
# the leading zero instead of coercing it to the number 1.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 993579326.0700001
$ws.Range("P2").Value = 58075286.05
$ws.Range("Q2").Value = 147147694.61
$ws.Range("R2").Value = -3.212520542
$ws.Range("S2").Value = 378460373.04
$ws.Range("T2").Value = 34.4166805516
$ws.Range("U2").Value = 105846690.34
$ws.Range("V2").Value = 17.2693486286
$ws.Range("W2").Value = 509768797.73
$ws.Range("X2").Value = 229905757.37
$ws.Range("Y2").Value = 21.2954485599
$ws.Range("Z2").Value = 16251974.19
$ws.Range("AA2").Value = 467.0117467074
$ws.Range("AB2").Value = 483810528.34
$ws.Range("AC2").Value = 14.0415905835
$ws.Range("AD2").Value = 21.4790108657
$ws.Range("AE2").Value = 29.4941540729
$ws.Range("AF2").Value = 146.0888744226
$ws.Range("AG2").Value = 51.306300801
